# Apply the "Issues fixes and reports" update:
#  - AMSIN  : append row 71 (173fnlrun sprint run)
#  - BETA   : append row 32 (173beta sprint run)
#  - AMS    : re-style existing row 35 (172live) to match the sheet's
#             normal data style, then append row 36 (live173 sprint run)

$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell without Excel's automatic
# date-recognition silently turning e.g. "2023-02-20" into a date serial
# number (which would also mint a brand new number-format style).
# We park the literal string behind a formula (so it is already a typed
# value, not "user keyboard input") in a scratch cell far outside the
# used range, copy it, and Paste-Special **values only** onto the real
# destination - this leaves the destination's existing number format /
# style completely untouched.
function Set-SafeText {
    param($ws, $row, $col, [string]$text)
    $scratch = $ws.Cells.Item(600, 30)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
    $scratch.Clear()
}

# Helper: append one new data row at $destRow by copying the formatting
# of $srcRow (Copy + Insert reproduces the exact existing style indices,
# unlike setting .Value/.NumberFormat which always mints new style
# records), then fills in the real values.
function Add-DataRow {
    param($ws, $srcRow, $destRow, $date, $time, $name, $total, $pass, $fail, $taken)

    $ws.Range("A$srcRow`:G$srcRow").Copy()
    $ws.Range("A$destRow`:G$destRow").Insert()

    Set-SafeText $ws $destRow 1 $date
    $ws.Cells.Item($destRow, 2).Value = $time
    $ws.Cells.Item($destRow, 3).Value = $name
    $ws.Cells.Item($destRow, 4).Value = $total
    $ws.Cells.Item($destRow, 5).Value = $pass
    $ws.Cells.Item($destRow, 6).Value = $fail
    $ws.Cells.Item($destRow, 7).Value = $taken
}

# ---------------------------------------------------------------------
# AMSIN: A1:G70 -> A1:G71
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")
Add-DataRow $wsAmsin 70 71 "2023-02-20" 44977.41628445602 "173fnlrun" 98 96 2 2.46

# ---------------------------------------------------------------------
# BETA: A1:G31 -> A1:G32
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
Add-DataRow $wsBeta 31 32 "2023-02-20" 44977.59682505787 "173beta" 98 98 0 2.55

# ---------------------------------------------------------------------
# AMS: A1:G35 -> A1:G36
#  Step 1: row 35 (172live) picks up the normal row styling (it was
#          previously unstyled) and its run-time value is refreshed.
#  Step 2: a fresh row 36 (live173) is appended.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Re-style row 35 in place: copy row 34's formatting in above row 35
# (pushing row 35's data down to row 36), rewrite row 35 with the
# original 172live values under the new formatting, then delete the
# now-duplicate shifted-down row.
$wsAms.Range("A34:G34").Copy()
$wsAms.Range("A35:G35").Insert()
Set-SafeText $wsAms 35 1 "2023-01-20"
$wsAms.Cells.Item(35, 2).Value = 44946.91931810185
$wsAms.Cells.Item(35, 3).Value = "172live"
$wsAms.Cells.Item(35, 4).Value = 98
$wsAms.Cells.Item(35, 5).Value = 96
$wsAms.Cells.Item(35, 6).Value = 2
$wsAms.Cells.Item(35, 7).Value = 3.71
$wsAms.Range("A36:G36").Delete()

# Append row 36 (live173). Only the Run Time cell (column B) carries the
# date-time number format in the source data, so only it is copied from
# an existing date-time cell; the rest are plain values.
$wsAms.Range("B35").Copy()
$wsAms.Range("B36").Insert()
Set-SafeText $wsAms 36 1 "2023-02-20"
$wsAms.Cells.Item(36, 2).Value = 44977.85398394011
$wsAms.Cells.Item(36, 3).Value = "live173"
$wsAms.Cells.Item(36, 4).Value = 98
$wsAms.Cells.Item(36, 5).Value = 98
$wsAms.Cells.Item(36, 6).Value = 0
$wsAms.Cells.Item(36, 7).Value = 2.68
